$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 507, shifting rows 507:540 down to 508:541
$ws.Rows.Item(507).Insert()

# Populate the new row 507 with its data
$ws.Cells.Item(507, 1).Value = 10
$ws.Cells.Item(507, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(507, 3).Value = "La Araucanía"
$ws.Cells.Item(507, 4).Value = 45021
$ws.Cells.Item(507, 4).NumberFormat = $ws.Cells.Item(508, 4).NumberFormat
$ws.Cells.Item(507, 5).Value = 9
$ws.Cells.Item(507, 6).Value = 100112040
$ws.Cells.Item(507, 7).Value = "Cilantro"
$ws.Cells.Item(507, 8).Value = "Sin especificar"
$ws.Cells.Item(507, 9).Value = "Primera"
$ws.Cells.Item(507, 10).Value = 55
$ws.Cells.Item(507, 11).Value = 7000
$ws.Cells.Item(507, 12).Value = 7000
$ws.Cells.Item(507, 13).Value = 7000
$ws.Cells.Item(507, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(507, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(507, 16).Value = 3500
$ws.Cells.Item(507, 17).Value = 2
$ws.Cells.Item(507, 18).Value = "Hortaliza"
